$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "- A couple years later, an expeditionary group"
#        -> "- A few years later, an expeditionary group"
#           (splitting the run around the replaced word, and placing the
#            "_GoBack" bookmark right after the newly typed word, mirroring
#            what Word itself produces when you select "couple" and type
#            "few")
# ---------------------------------------------------------------------

# Locate the run of text that contains the phrase we need to edit.
$full = $d.Content
$full.Find.Execute("- A couple years later, an expeditionary group") | Out-Null
$phraseStart = $full.Start
$phraseEnd = $full.End

# Find "couple" within that phrase only, so we don't disturb anything else.
$wordRange = $d.Range($phraseStart, $phraseEnd)
$wordRange.Find.Execute("couple") | Out-Null
$wordStart = $wordRange.Start
$wordEnd = $wordRange.End

# Drop zero-length bookmarks right before/after "couple" - inserting a
# bookmark forces the surrounding run to split at that exact point without
# touching formatting, which is how we get three separate runs instead of
# one merged run.
$d.Bookmarks.Add("TempSplitBefore", $d.Range($wordStart, $wordStart)) | Out-Null
$d.Bookmarks.Add("TempSplitAfter", $d.Range($wordEnd, $wordEnd)) | Out-Null

# Replace just "couple" with "few" (range is now tightly bounded by the two
# bookmarks we just added, so the edit cannot spill into neighboring runs).
$midRange = $d.Range($wordStart, $wordEnd)
$midRange.Text = "few"

# Drop the helper bookmark before the word - it was only needed to force
# the left-hand split.
$d.Bookmarks("TempSplitBefore").Delete()

# Promote the helper bookmark after the word into the real "_GoBack"
# bookmark (Word keeps only a single "_GoBack" bookmark, so adding one
# under that name automatically relocates/removes any existing one,
# including the one that used to sit near "Nordic Inn").
$afterRange = $d.Bookmarks("TempSplitAfter").Range
$d.Bookmarks.Add("_GoBack", $afterRange) | Out-Null
$d.Bookmarks("TempSplitAfter").Delete()
